$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1722054380664653
$ws.Range("C2").Value = 0.5891238670694864
$ws.Range("J2").Value = 0.006042296072507553
$ws.Range("P2").Value = 0.1359516616314199
$ws.Range("S2").Value = 0.09667673716012085
$ws.Range("C3").Value = 0.01020408163265306
$ws.Range("J3").Value = 0.03061224489795918
$ws.Range("P3").Value = 0.7295918367346939
$ws.Range("S3").Value = 0.2295918367346939
$ws.Range("J4").Value = 0.05405405405405406
$ws.Range("P4").Value = 0.7567567567567568
$ws.Range("S4").Value = 0.1891891891891892
$ws.Range("B6").Value = 0.08294930875576037
$ws.Range("D6").Value = 0.0184331797235023
$ws.Range("F6").Value = 0.06912442396313365
$ws.Range("J6").Value = 0.2350230414746544
$ws.Range("O6").Value = 0.01382488479262673
$ws.Range("Q6").Value = 0.1566820276497696
$ws.Range("R6").Value = 0.05990783410138249
$ws.Range("S6").Value = 0.3640552995391705
$ws.Range("B7").Value = 0.09615384615384616
$ws.Range("D7").Value = 0.01442307692307692
$ws.Range("F7").Value = 0.0673076923076923
$ws.Range("J7").Value = 0.1105769230769231
$ws.Range("O7").Value = 0.01923076923076923
$ws.Range("Q7").Value = 0.1971153846153846
$ws.Range("R7").Value = 0.08173076923076923
$ws.Range("S7").Value = 0.4134615384615384
$ws.Range("B8").Value = 0.1233140655105973
$ws.Range("D8").Value = 0.005780346820809248
$ws.Range("F8").Value = 0.06936416184971098
$ws.Range("J8").Value = 0.1136801541425819
$ws.Range("O8").Value = 0.01734104046242774
$ws.Range("Q8").Value = 0.1695568400770713
$ws.Range("R8").Value = 0.0905587668593449
$ws.Range("S8").Value = 0.4104046242774567
$ws.Range("B9").Value = 0.131578947368421
$ws.Range("D9").Value = 0.01754385964912281
$ws.Range("F9").Value = 0.08771929824561403
$ws.Range("J9").Value = 0.06140350877192982
$ws.Range("O9").Value = 0.03508771929824561
$ws.Range("Q9").Value = 0.1140350877192982
$ws.Range("R9").Value = 0.1140350877192982
$ws.Range("S9").Value = 0.4385964912280702
$ws.Range("B10").Value = 0.1164225134926754
$ws.Range("D10").Value = 0.02235929067077872
$ws.Range("F10").Value = 0.06168080185042406
$ws.Range("J10").Value = 0.1441788743253662
$ws.Range("O10").Value = 0.01002313030069391
$ws.Range("Q10").Value = 0.2259059367771781
$ws.Range("R10").Value = 0.06245181187355436
$ws.Range("S10").Value = 0.3569776407093292
$ws.Range("G11").Value = 0.1666666666666667
$ws.Range("J11").Value = 0.1121212121212121
$ws.Range("K11").Value = 0.2393939393939394
$ws.Range("L11").Value = 0.4727272727272727
$ws.Range("S11").Value = 0.00909090909090909
$ws.Range("G12").Value = 0.7891566265060241
$ws.Range("J12").Value = 0.1325301204819277
$ws.Range("K12").Value = 0.006024096385542169
$ws.Range("L12").Value = 0.04216867469879518
$ws.Range("S12").Value = 0.03012048192771084
$ws.Range("G13").Value = 0.4583333333333333
$ws.Range("J13").Value = 0.4166666666666667
$ws.Range("S13").Value = 0.125
$ws.Range("F15").Value = 0.004901960784313725
$ws.Range("H15").Value = 0.196078431372549
$ws.Range("I15").Value = 0.01470588235294118
$ws.Range("J15").Value = 0.392156862745098
$ws.Range("K15").Value = 0.07352941176470588
$ws.Range("M15").Value = 0.009803921568627451
$ws.Range("O15").Value = 0.03431372549019608
$ws.Range("S15").Value = 0.2745098039215687
$ws.Range("F16").Value = 0.01923076923076923
$ws.Range("H16").Value = 0.1634615384615385
$ws.Range("J16").Value = 0.4278846153846154
$ws.Range("K16").Value = 0.1009615384615385
$ws.Range("M16").Value = 0.02884615384615385
$ws.Range("N16").Value = 0.004807692307692308
$ws.Range("O16").Value = 0.02884615384615385
$ws.Range("S16").Value = 0.1634615384615385
$ws.Range("F17").Value = 0.02155172413793104
$ws.Range("H17").Value = 0.2155172413793103
$ws.Range("I17").Value = 0.05172413793103448
$ws.Range("J17").Value = 0.4245689655172414
$ws.Range("K17").Value = 0.09267241379310345
$ws.Range("M17").Value = 0.01724137931034483
$ws.Range("N17").Value = 0.002155172413793103
$ws.Range("O17").Value = 0.07758620689655173
$ws.Range("S17").Value = 0.09698275862068965
$ws.Range("F18").Value = 0.01796407185628742
$ws.Range("H18").Value = 0.1856287425149701
$ws.Range("I18").Value = 0.0658682634730539
$ws.Range("J18").Value = 0.4491017964071856
$ws.Range("K18").Value = 0.08383233532934131
$ws.Range("M18").Value = 0.02395209580838323
$ws.Range("O18").Value = 0.0658682634730539
$ws.Range("S18").Value = 0.1077844311377246
$ws.Range("F19").Value = 0.0132398753894081
$ws.Range("H19").Value = 0.2406542056074766
$ws.Range("I19").Value = 0.04906542056074766
$ws.Range("J19").Value = 0.3637071651090343
$ws.Range("K19").Value = 0.1214953271028037
$ws.Range("M19").Value = 0.02180685358255452
$ws.Range("N19").Value = 0.002336448598130841
$ws.Range("O19").Value = 0.06386292834890965
$ws.Range("S19").Value = 0.1238317757009346
